# Summary of the edit (per the target diff / commit message "added better
# dependency tracking"):
#   - Sheet3!A1 gets a new cross-sheet formula that depends on two cells of
#     Sheet1 (=Sheet1!A2+Sheet1!E10, evaluates to 5).
#   - Sheet3 becomes the active sheet/tab (was Sheet1).
#   - The selection on Sheet1 moves from E4 to E10 (the cell the new formula
#     references last).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Move to Sheet1 and land the selection on E10 (one of the precedents of the
# formula we're about to write).
$ws1.Activate()
$ws1.Range("E10").Select()

# Switch to Sheet3 and enter the new dependency-tracking formula in A1.
$ws3.Activate()
$ws3.Range("A1").Formula = "=Sheet1!A2+Sheet1!E10"
$ws3.Range("A1").Select()
